# Updates the feature_importance worksheet (Sheet1) with the refreshed
# 2024-season (weeks 1-4) feature importances and reordered feature list,
# per the commit: "updated dataset for 2024 season 1-4 / also created new
# version of Model without some of the stats".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; A='Home_Last 5_Margin'; B=0.05351933464407921},
    @{Row=3; A='Away _Last 5_ Margin'; B=0.04466285929083824},
    @{Row=4; A='A_Prev_PPG'; B=0.04288320243358612},
    @{Row=5; A='A_Prev_Tot_Def'; B=0.04272532090544701},
    @{Row=6; A='A_Prev_Pass'; B=0.0414387471973896},
    @{Row=7; A='Prev_Sacks_Taken'; B=0.03404835984110832},
    @{Row=8; A='Def_PassRTG'; B=0.03227928280830383},
    @{Row=9; A='A_Prev_Sacks_Taken'; B=0.03184231370687485},
    @{Row=10; A='Prev_Rush'; B=0.03159726411104202},
    @{Row=11; A='A_Prev_Rush'; B=0.03113587573170662},
    @{Row=12; A='AwayDiv_AFC West'; B=0.03051367402076721},
    @{Row=13; A='Prev_PPG'; B=0.02978039532899857},
    @{Row=14; A='HomeDiv_AFC North'; B=0.02932943776249886},
    @{Row=15; A='A_Prev_Def_PPG'; B=0.02738093212246895},
    @{Row=16; A='Away_PrevYrW'; B=0.02699758298695087},
    @{Row=17; A='Home_QBInjury'; B=0.02596757374703884},
    @{Row=18; A='Home_PrevYrW'; B=0.0257403627038002},
    @{Row=19; A='Away_QBInjury'; B=0.02519631385803223},
    @{Row=20; A='Week_WinInteraction'; B=0.02432323805987835},
    @{Row=21; A='Away_Third_Down'; B=0.0237570833414793},
    @{Row=22; A='Away_WinInteraction'; B=0.02296528592705727},
    @{Row=23; A='AwayDiv_NFC West'; B=0.02284062094986439},
    @{Row=24; A='A_Prev_Spec'; B=0.02257499098777771},
    @{Row=25; A='Prev_Def_PPG'; B=0.0223784688860178},
    @{Row=26; A='Week_TO_Margin_Interaction'; B=0.02150846086442471},
    @{Row=27; A='Prev_Tot_Def'; B=0.02126043289899826},
    @{Row=28; A='Week_Third_Down'; B=0.0208237674087286},
    @{Row=29; A='HomeDiv_NFC North'; B=0.02044583857059479},
    @{Row=30; A='A_Wins5'; B=0.01979056186974049},
    @{Row=31; A='Day_Mon'; B=0.01970954798161983},
    @{Row=32; A='Away_TO_Margin_Interaction'; B=0.01957257837057114},
    @{Row=33; A='HomeDiv_NFC East'; B=0.01946909911930561},
    @{Row=34; A='Away_Penalty_Yards'; B=0.01933771558105946},
    @{Row=35; A='Prev_Pass'; B=0.01843366026878357},
    @{Row=36; A='A_Def_PassRTG'; B=0.01728787273168564},
    @{Row=37; A='Prev_Spec'; B=0.01265790220350027},
    @{Row=38; A='Day_Thu'; B=0.01179203484207392},
    @{Row=39; A='HomeDiv_NFC West'; B=0.008180109784007072},
    @{Row=40; A='Week'; B=0.003851997898891568},
    @{Row=41; A='AwayDiv_NFC South'; B=0},
    @{Row=42; A='HomeDiv_AFC West'; B=0},
    @{Row=43; A='HomeDiv_AFC South'; B=0},
    @{Row=44; A='Day_Sat'; B=0},
    @{Row=45; A='AwayDiv_AFC North'; B=0},
    @{Row=46; A='AwayDiv_NFC North'; B=0},
    @{Row=47; A='Day_Tue'; B=0},
    @{Row=48; A='Day_Wed'; B=0},
    @{Row=49; A='HomeTeam'; B=0},
    @{Row=50; A='AwayDiv_NFC East'; B=0},
    @{Row=51; A='HomeDiv_NFC South'; B=0},
    @{Row=52; A='AwayDiv_AFC South'; B=0},
    @{Row=53; A='H_Wins5'; B=0},
    @{Row=54; A='AwayTeam'; B=0},
    @{Row=55; A='Day_Sun'; B=0}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 1).Value = $u.A
    $ws.Cells.Item($u.Row, 2).Value = $u.B
}
